# Executive presentation trim:
# Keep only the Title slide, "Why This Solution?", "Business Value - Financial
# Impact" and "Risk Mitigation" slides; drop every other slide. Also strip the
# now-unwanted bold emphasis from a handful of table cells on the three
# surviving content slides.

$p = $ppt.ActivePresentation

# --- 1. Remove slides we no longer want -------------------------------------------------
# Original deck order (1-based):
#   1  Title                                   -> keep
#   2  Agenda                                   -> remove
#   3  Executive Summary                        -> remove
#   4  Current State - Business Challenge       -> remove
#   5  Vision - Future State                    -> remove
#   6  Solution Overview                        -> remove
#   7  Why This Solution?                       -> keep (becomes slide 2)
#   8  Business Value - Financial Impact        -> keep (becomes slide 3)
#   9  Business Value - Strategic Benefits      -> remove
#   10 Implementation Approach                  -> remove
#   11 Risk Mitigation                          -> keep (becomes slide 4)
#   12 Investment Summary                       -> remove
#   13 Timeline & Milestones                    -> remove
#   14 Success Stories                          -> remove
#   15 Our Partnership Advantage                -> remove
#   16 Next Steps                               -> remove
#   17 Contact Information                      -> remove
#
# Delete from the highest index down so earlier indices stay valid.
$slidesToRemove = @(17, 16, 15, 14, 13, 12, 10, 9, 6, 5, 4, 3, 2)
foreach ($idx in $slidesToRemove) {
    $p.Slides.Item($idx).Delete()
}

# --- 2. Strip bold from specific table cells --------------------------------------------
# Slide 2 ("Why This Solution?"): the bottom row ([Current limitation 3] /
# checkmarked [Our advantage 3]) should no longer be bold.
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
$tbl2.Cell(4, 1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl2.Cell(4, 2).Shape.TextFrame.TextRange.Font.Bold = $false

# Slide 3 ("Business Value - Financial Impact"): header row (Metric/Value) and
# the ROI row lose their bold emphasis.
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
$tbl3.Cell(1, 1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell(1, 2).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell(6, 1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl3.Cell(6, 2).Shape.TextFrame.TextRange.Font.Bold = $false

# Slide 4 ("Risk Mitigation"): header row (Risk/Mitigation Strategy/Success
# Probability) and the [Risk 3] row lose their bold emphasis.
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
$tbl4.Cell(1, 1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(1, 2).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(1, 3).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(4, 1).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(4, 2).Shape.TextFrame.TextRange.Font.Bold = $false
$tbl4.Cell(4, 3).Shape.TextFrame.TextRange.Font.Bold = $false

Write-Host "Final slide count: $($p.Slides.Count)"
